# Weekly CompStat refresh (105th Precinct) - new crime data collected.
# Updates: report header (volume number + covering-week dates) and the
# weekly/28-day/YTD/2yr/15yr/32yr crime-category figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: "Volume 32   Number  30" -> "...31"
# ---------------------------------------------------------------------
$ws.Range("A8").Characters(21, 2).Text = "31"
$aPart1 = $ws.Range("A8").Characters(1, 21)
$aPart1.Font.Size = 10
$aPart1.Font.Name = "Andale WT"
$aPart2 = $ws.Range("A8").Characters(22, 1)
$aPart2.Font.Size = 10
$aPart2.Font.Name = "Andale WT"

# ---------------------------------------------------------------------
# Header text: "Report Covering the Week  7/21/2025  Through  7/27/2025"
#           -> "Report Covering the Week  7/28/2025  Through  8/3/2025"
# ---------------------------------------------------------------------
$ws.Range("C9").Characters(27, 9).Text = "7/28/2025"
$ws.Range("C9").Characters(47, 9).Text = "8/3/2025"
$cPart1 = $ws.Range("C9").Characters(1, 53)
$cPart1.Font.Size = 10
$cPart1.Font.Name = "Andale WT"
$cPart2 = $ws.Range("C9").Characters(54, 1)
$cPart2.Font.Size = 10
$cPart2.Font.Name = "Andale WT"

# ---------------------------------------------------------------------
# Crime-complaint figures
# ---------------------------------------------------------------------

# Row 15 - Rape
$ws.Range("F15").Value = 4
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 11
$ws.Range("K15").Value = -15.384615384615
$ws.Range("L15").Value = 37.5
$ws.Range("M15").Value = -8.333333333333
$ws.Range("N15").Value = -56

# Row 16 - Robbery
$ws.Range("F16").Value = 3
$ws.Range("H16").Value = -40
$ws.Range("J16").Value = 30
$ws.Range("K16").Value = 53.333333333333
$ws.Range("L16").Value = -16.363636363636
$ws.Range("M16").Value = -76.884422110552
$ws.Range("N16").Value = -92.396694214876

# Row 17 - Fel. Assault
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = 21.428571428571
$ws.Range("I17").Value = 124
$ws.Range("J17").Value = 156
$ws.Range("K17").Value = -20.512820512820
$ws.Range("L17").Value = -16.216216216216
$ws.Range("M17").Value = -31.868131868131
$ws.Range("N17").Value = -46.551724137931

# Row 18 - Burglary
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = -42.857142857142
$ws.Range("I18").Value = 82
$ws.Range("J18").Value = 69
$ws.Range("K18").Value = 18.840579710144
$ws.Range("L18").Value = -4.651162790697
$ws.Range("M18").Value = -60.952380952380
$ws.Range("N18").Value = -91.211146838156

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 11.111111111111
$ws.Range("F19").Value = 26
$ws.Range("G19").Value = 27
$ws.Range("H19").Value = -3.703703703703
$ws.Range("I19").Value = 212
$ws.Range("J19").Value = 219
$ws.Range("K19").Value = -3.196347031963
$ws.Range("L19").Value = -15.537848605577
$ws.Range("M19").Value = -17.509727626459
$ws.Range("N19").Value = -40.112994350282

# Row 20 - G.L.A.
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -60
$ws.Range("F20").Value = 14
$ws.Range("G20").Value = 22
$ws.Range("H20").Value = -36.363636363636
$ws.Range("I20").Value = 109
$ws.Range("J20").Value = 142
$ws.Range("K20").Value = -23.239436619718
$ws.Range("L20").Value = 15.957446808510
$ws.Range("M20").Value = -51.555555555555
$ws.Range("N20").Value = -94.441611422743

# Row 21 - TOTAL
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = 4.761904761904
$ws.Range("F21").Value = 72
$ws.Range("G21").Value = 84
$ws.Range("H21").Value = -14.285714285714
$ws.Range("I21").Value = 585
$ws.Range("J21").Value = 629
$ws.Range("K21").Value = -6.995230524642
$ws.Range("L21").Value = -9.020217729393
$ws.Range("M21").Value = -46.526508226691
$ws.Range("N21").Value = -85.807860262008

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 6
$ws.Range("D24").Value = 13
$ws.Range("E24").Value = -53.846153846153
$ws.Range("F24").Value = 44
$ws.Range("G24").Value = 51
$ws.Range("H24").Value = -13.725490196078
$ws.Range("I24").Value = 344
$ws.Range("J24").Value = 412
$ws.Range("K24").Value = -16.504854368932
$ws.Range("L24").Value = -26.652452025586
$ws.Range("M24").Value = -29.218106995884

# Row 25 - Retail Theft
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 14
$ws.Range("G25").Value = 10
$ws.Range("H25").Value = 40
$ws.Range("I25").Value = 91
$ws.Range("J25").Value = 102
$ws.Range("K25").Value = -10.784313725490
$ws.Range("L25").Value = 9.638554216867

# Row 26 - Misd. Assault
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 100
$ws.Range("F26").Value = 26
$ws.Range("G26").Value = 29
$ws.Range("H26").Value = -10.344827586206
$ws.Range("I26").Value = 251
$ws.Range("J26").Value = 251
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 21.256038647343
$ws.Range("M26").Value = -35.309278350515

# Row 27 - UCR Rape*
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 11
$ws.Range("K27").Value = -38.888888888888
$ws.Range("L27").Value = -8.333333333333

# Row 28 - Other Sex Crimes
$ws.Range("F28").Value = 1
$ws.Range("H28").Value = 0
$ws.Range("L28").Value = 42.857142857142

# ---------------------------------------------------------------------
# Row 29 - Shooting Vic.: 28-Day figures now show no data ("0" / "***.*")
# ---------------------------------------------------------------------
$ws.Range("I29").Copy() | Out-Null
$ws.Range("G29").PasteSpecial(-4122) | Out-Null
$ws.Range("I29").Copy() | Out-Null
$ws.Range("G29").PasteSpecial(-4163) | Out-Null

$ws.Range("E29").Copy() | Out-Null
$ws.Range("H29").PasteSpecial(-4122) | Out-Null
$ws.Range("E29").Copy() | Out-Null
$ws.Range("H29").PasteSpecial(-4163) | Out-Null

# ---------------------------------------------------------------------
# Row 30 - Shooting Inc.: same change as row 29
# ---------------------------------------------------------------------
$ws.Range("I30").Copy() | Out-Null
$ws.Range("G30").PasteSpecial(-4122) | Out-Null
$ws.Range("I30").Copy() | Out-Null
$ws.Range("G30").PasteSpecial(-4163) | Out-Null

$ws.Range("E30").Copy() | Out-Null
$ws.Range("H30").PasteSpecial(-4122) | Out-Null
$ws.Range("E30").Copy() | Out-Null
$ws.Range("H30").PasteSpecial(-4163) | Out-Null

# ---------------------------------------------------------------------
# Row 33 - Traffic Fatalities: week-to-date 28-Day 2025 figure now "0"
# ---------------------------------------------------------------------
$ws.Range("G33").Copy() | Out-Null
$ws.Range("F33").PasteSpecial(-4122) | Out-Null
$ws.Range("G33").Copy() | Out-Null
$ws.Range("F33").PasteSpecial(-4163) | Out-Null
